$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# The sheet had duplicate/extra credential rows (old rows 5-7). Consolidate
# them into a single row 5 (Vendoadmin / Vendo@2022), matching what used to
# be row 7, by removing old row 6 (A6 only) and old row 5 (B5 only) - this
# shifts old row 7's data+style up into row 5 untouched.
$ws.Rows(6).Delete()
$ws.Rows(5).Delete()

# Clean up the now-stale hyperlink entry left pointing at the deleted B7.
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$7') {
        $hl.Delete()
    }
}

# Preserve the sheet's prior max row-outline level (6 -> 4) in
# sheetFormatPr without leaving any visible outlineLevel on real rows.
$ws.Rows(100).OutlineLevel = 4
$ws.Rows(100).Delete()

$ws.Range("A1").Select()
